$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift existing columns to make room for new data ---
# Insert 6 new columns before A: old A..I -> G..O ; old L..P -> R..V
$ws.Columns("A:F").Insert()
# Insert 4 more columns before (current) P so the "Failure" helper block
# (currently at R..V) lands at V..Z
$ws.Columns("P:S").Insert()

# --- 2. New header row (row 1) ---
$ws.Range("A1").Value = "casing_flow"
$ws.Range("B1").Value = "tubing_flow"
$ws.Range("C1").Value = "d_production_casing_in"
$ws.Range("D1").Value = "d_tubing_in"
$ws.Range("E1").Value = "d_production_casing"
$ws.Range("F1").Value = "d_tubing"
$ws.Range("G1").Value = "mode"
$ws.Range("H1").Value = "moment_crit_rup_conductor"
$ws.Range("J1").Value = "moment_crit_rup_production"
$ws.Range("K1").Value = "moment_crit_rup_tubing"
$ws.Range("L1").Value = "sigma_moment_crit_rup_conductor"
$ws.Range("M1").Value = "sigma_moment_crit_rup_surface"
$ws.Range("N1").Value = "sigma_moment_crit_rup_production"
$ws.Range("O1").Value = "sigma_moment_crit_rup_tubing"
$ws.Range("P1").Value = "sigma_mu_moment_crit_rup_conductor"
$ws.Range("Q1").Value = "sigma_mu_moment_crit_rup_surface"
$ws.Range("R1").Value = "sigma_mu_moment_crit_rup_production"
$ws.Range("S1").Value = "sigma_mu_moment_crit_rup_tubing"
$ws.Range("I1").Value = "moment_crit_rup_surface"

$ws.Range("A1:G1").WrapText = $true
$ws.Rows("1:1").RowHeight = 45

# --- 3. New input columns A-F, rows 2-4 ---
$ws.Range("A2").Value = $true
$ws.Range("A3").Value = $true
$ws.Range("A4").Value = $true

$ws.Range("B2").Formula = '=IF(A2=TRUE,FALSE,TRUE)'
$ws.Range("B3:B4").Formula = '=IF(A3=TRUE,FALSE,TRUE)'

$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 4
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 3.5
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 3

$ws.Range("E2").Formula = '=C2*2.54/100'
$ws.Range("F2").Formula = '=D2*2.54/100'
$ws.Range("E3:F4").Formula = '=C3*2.54/100'

$ws.Range("G2:G4").Formula = '=IF(B2=TRUE,IF(D2<(3+1/8),4,IF(C2>(7+3/4),1,2)),IF(C2>=(8+5/8),1,IF(C2<=(6+5/8),4,2)))'

# --- 4. Shifted "moment_crit_rup_*" lookup formulas (H:K) ---
$ws.Range("H2:K4").Formula = '=IF($G2=1,W$3,IF($G2=2,W$4,IF($G2=4,W$5)))'

# --- 5. New sigma_mu columns P:S ---
$ws.Range("P2:S4").Value = 0.25

# --- 6. Formatting / view metadata ---
$ws.Range("A1").Select()
$ws.Range("A1:S4").Select()

$wb.Application.ActiveWindow.ScrollColumn = 3
